$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target table (rows 2-26), columns A-T: a new "Resolving-Mac" sending-cluster
# block (rows 22-26) was added and TPM-derived values were recomputed throughout.
$data = @(
    @("ECs", "Lamc2", "Itgb4", "ECs", 3, 1, 1.598137666666667, 4.794413, 0.2574006098591189, 0.2574006098591188, 3, 1, 15.991592, 47.97477600000001, 0.7512457858594542, 0.7512457858594542, 25.55676552516534, 230.010889726488, 0.1933711234343166, 0.1933711234343165),
    @("ECs", "Lamc2", "Itgb4", "FAPs", 3, 1, 1.598137666666667, 4.794413, 0.2574006098591189, 0.2574006098591188, 3, 1, 2.388815333333333, 7.166446, 0.1122206877441041, 0.1122206877441042, 3.817655762910889, 34.358901866198, 0.02888567346414216, 0.02888567346414215),
    @("ECs", "Lamc2", "Itgb4", "Inflammatory-Mac", 3, 1, 1.598137666666667, 4.794413, 0.2574006098591189, 0.2574006098591188, 3, 1, 0.914392, 2.743176, 0.04295589436146182, 0.04295589436146183, 1.461324297298667, 13.151918675688, 0.01105687340568416, 0.01105687340568416),
    @("ECs", "Lamc2", "Itgb4", "MuSCs", 3, 1, 1.598137666666667, 4.794413, 0.2574006098591189, 0.2574006098591188, 3, 1, 1.309904666666666, 3.929714, 0.06153610977011958, 0.06153610977011959, 2.093407987542444, 18.840671887882, 0.01583943218318647, 0.01583943218318646),
    @("ECs", "Lamc2", "Itgb4", "Resolving-Mac", 3, 1, 1.598137666666667, 4.794413, 0.2574006098591189, 0.2574006098591188, 3, 1, 0.6820603333333333, 2.046181, 0.03204152226486026, 0.03204152226486026, 1.090026309639222, 9.810236786753, 0.008247507371789567, 0.008247507371789566),
    @("FAPs", "Lamc2", "Itgb4", "ECs", 3, 1, 3.908618333333334, 11.725855, 0.6295332146228534, 0.6295332146228533, 3, 1, 15.991592, 47.97477600000001, 0.7512457858594542, 0.7512457858594542, 62.50502967038668, 562.5452670334801, 0.472934174543974, 0.4729341745439739),
    @("FAPs", "Lamc2", "Itgb4", "FAPs", 3, 1, 3.908618333333334, 11.725855, 0.6295332146228534, 0.6295332146228533, 3, 1, 2.388815333333333, 7.166446, 0.1122206877441041, 0.1122206877441042, 9.336967406814443, 84.03270666133, 0.07064665030273333, 0.07064665030273332),
    @("FAPs", "Lamc2", "Itgb4", "Inflammatory-Mac", 3, 1, 3.908618333333334, 11.725855, 0.6295332146228534, 0.6295332146228533, 3, 1, 0.914392, 2.743176, 0.04295589436146182, 0.04295589436146183, 3.574009335053333, 32.16608401548, 0.02704216226437076, 0.02704216226437076),
    @("FAPs", "Lamc2", "Itgb4", "MuSCs", 3, 1, 3.908618333333334, 11.725855, 0.6295332146228534, 0.6295332146228533, 3, 1, 1.309904666666666, 3.929714, 0.06153610977011958, 0.06153610977011959, 5.119917395052222, 46.07925655547, 0.03873902499896816, 0.03873902499896816),
    @("FAPs", "Lamc2", "Itgb4", "Resolving-Mac", 3, 1, 3.908618333333334, 11.725855, 0.6295332146228534, 0.6295332146228533, 3, 1, 0.6820603333333333, 2.046181, 0.03204152226486026, 0.03204152226486026, 2.665913523306111, 23.993221709755, 0.02017120251280721, 0.02017120251280721),
    @("Inflammatory-Mac", "Lamc2", "Itgb4", "ECs", 3, 1, 0.1524303333333333, 0.457291, 0.02455086415857089, 0.02455086415857088, 3, 1, 15.991592, 47.97477600000001, 0.7512457858594542, 0.7512457858594542, 2.437603699090667, 21.938433291816, 0.01844373323833429, 0.01844373323833429),
    @("Inflammatory-Mac", "Lamc2", "Itgb4", "FAPs", 3, 1, 0.1524303333333333, 0.457291, 0.02455086415857089, 0.02455086415857088, 3, 1, 2.388815333333333, 7.166446, 0.1122206877441041, 0.1122206877441042, 0.3641279175317778, 3.277151257786, 0.002755114860586902, 0.002755114860586901),
    @("Inflammatory-Mac", "Lamc2", "Itgb4", "Inflammatory-Mac", 3, 1, 0.1524303333333333, 0.457291, 0.02455086415857089, 0.02455086415857088, 3, 1, 0.914392, 2.743176, 0.04295589436146182, 0.04295589436146183, 0.1393810773573333, 1.254429696216, 0.00105460432727817, 0.00105460432727817),
    @("Inflammatory-Mac", "Lamc2", "Itgb4", "MuSCs", 3, 1, 0.1524303333333333, 0.457291, 0.02455086415857089, 0.02455086415857088, 3, 1, 1.309904666666666, 3.929714, 0.06153610977011958, 0.06153610977011959, 0.1996692049748889, 1.797022844774, 0.001510764671813113, 0.001510764671813112),
    @("Inflammatory-Mac", "Lamc2", "Itgb4", "Resolving-Mac", 3, 1, 0.1524303333333333, 0.457291, 0.02455086415857089, 0.02455086415857088, 3, 1, 0.6820603333333333, 2.046181, 0.03204152226486026, 0.03204152226486026, 0.1039666839634444, 0.935700155671, 0.0007866470605584089, 0.0007866470605584087),
    @("MuSCs", "Lamc2", "Itgb4", "ECs", 3, 1, 0.5484213333333333, 1.645264, 0.0883303038305739, 0.08833030383057389, 3, 1, 15.991592, 47.97477600000001, 0.7512457858594542, 0.7512457858594542, 8.770130206762667, 78.93117186086401, 0.06635776851640385, 0.06635776851640383),
    @("MuSCs", "Lamc2", "Itgb4", "FAPs", 3, 1, 0.5484213333333333, 1.645264, 0.0883303038305739, 0.08833030383057389, 3, 1, 2.388815333333333, 7.166446, 0.1122206877441041, 0.1122206877441042, 1.310077290193778, 11.790695611744, 0.00991248744451268, 0.00991248744451268),
    @("MuSCs", "Lamc2", "Itgb4", "Inflammatory-Mac", 3, 1, 0.5484213333333333, 1.645264, 0.0883303038305739, 0.08833030383057389, 3, 1, 0.914392, 2.743176, 0.04295589436146182, 0.04295589436146183, 0.5014720798293333, 4.513248718464, 0.003794307200261959, 0.003794307200261959),
    @("MuSCs", "Lamc2", "Itgb4", "MuSCs", 3, 1, 0.5484213333333333, 1.645264, 0.0883303038305739, 0.08833030383057389, 3, 1, 1.309904666666666, 3.929714, 0.06153610977011958, 0.06153610977011959, 0.7183796638328888, 6.465416974496, 0.00543550327254621, 0.00543550327254621),
    @("MuSCs", "Lamc2", "Itgb4", "Resolving-Mac", 3, 1, 0.5484213333333333, 1.645264, 0.0883303038305739, 0.08833030383057389, 3, 1, 0.6820603333333333, 2.046181, 0.03204152226486026, 0.03204152226486026, 0.3740564374204444, 3.366507936784, 0.002830237396849205, 0.002830237396849205),
    @("Resolving-Mac", "Lamc2", "Itgb4", "ECs", 1, 0.3333333333333333, 0.001148666666666667, 0.003446, 0.0001850075288829985, 0.0001850075288829985, 3, 1, 15.991592, 47.97477600000001, 0.7512457858594542, 0.7512457858594542, 0.01836900867733333, 0.165321078096, 0.0001389861264256239, 0.0001389861264256239),
    @("Resolving-Mac", "Lamc2", "Itgb4", "FAPs", 1, 0.3333333333333333, 0.001148666666666667, 0.003446, 0.0001850075288829985, 0.0001850075288829985, 3, 1, 2.388815333333333, 7.166446, 0.1122206877441041, 0.1122206877441042, 0.002743952546222222, 0.024695572916, 0.0000207616721290873, 0.0000207616721290873),
    @("Resolving-Mac", "Lamc2", "Itgb4", "Inflammatory-Mac", 1, 0.3333333333333333, 0.001148666666666667, 0.003446, 0.0001850075288829985, 0.0001850075288829985, 3, 1, 0.914392, 2.743176, 0.04295589436146182, 0.04295589436146183, 0.001050331610666667, 0.009452984495999999, 0.00000794716386677318, 0.00000794716386677318),
    @("Resolving-Mac", "Lamc2", "Itgb4", "MuSCs", 1, 0.3333333333333333, 0.001148666666666667, 0.003446, 0.0001850075288829985, 0.0001850075288829985, 3, 1, 1.309904666666666, 3.929714, 0.06153610977011958, 0.06153610977011959, 0.001504643827111111, 0.013541794444, 0.00001138464360564276, 0.00001138464360564277),
    @("Resolving-Mac", "Lamc2", "Itgb4", "Resolving-Mac", 1, 0.3333333333333333, 0.001148666666666667, 0.003446, 0.0001850075288829985, 0.0001850075288829985, 3, 1, 0.6820603333333333, 2.046181, 0.03204152226486026, 0.03204152226486026, 0.0007834599695555554, 0.007051139725999999, 0.000005927922855871374, 0.000005927922855871374)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $c + 1).Value = $row[$c]
    }
}
